# Refresh the crypto price/volume snapshot (GitHub Actions style update).
# Column D ("Price") cells are leading-apostrophe-quoted so Excel keeps them
# as literal text (matching the source data) instead of silently coercing
# numeric-looking strings to Number and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.638.74"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "'2.167.92"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'226.76"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'63.14"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").Value = "'2.490.70"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "'21.78"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "'2.166.93"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "'39.619.06"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "'0.0₃0914"
$ws.Range("E19").Value = "  +7.01%  "
$ws.Range("D20").Value = "'71.78"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'6.02"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.34"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "'170.80"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'19.70"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("E34").Value = "  -2.83%  "
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").Value = "'3.80"
$ws.Range("E37").Value = "  +5.82%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'4.91"
$ws.Range("E40").Value = "  +17.72%  "
$ws.Range("D41").Value = "'102.44"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").Value = "'17.67"
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("D44").Value = "'1.514.74"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'0.0920"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").Value = "'0.000196"
$ws.Range("E50").Value = "  +32.14%  "
$ws.Range("D51").Value = "'2.369.33"
$ws.Range("E51").Value = "  +0.34%  "
